# Modify and Improve Generate Report Logic
# Update the single data row (row 2) of the assessment-number report:
#   cp_id:            109 -> 111
#   client_name:       Anantara_Analytics -> Hogwarts
#   project_name:      NeuCode -> Gryffindoar
#   assessment_type:   "360" -> "180"  (kept as text, same as the source data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111
$ws.Range("B2").Value = "Hogwarts"
$ws.Range("C2").Value = "Gryffindoar"

# assessment_type must stay a text value ("180"), not become a number.
# Force text formatting for the write, then clear the formatting change
# back off so we don't leave a stray number-format on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "180"
$ws.Range("D2").ClearFormats()
